$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: wipe all cell contents (keeps per-cell formatting/styles
# intact) so the shared-strings table gets rebuilt from scratch as we
# re-enter the values below, in the same row-major / column-major order
# the corrected import file was produced in.
$ws.Cells.ClearContents()

# ---- Row 1 (headers) ----
$ws.Range("A1").Value = "EAN"
$ws.Range("B1").Value = "NOMBRE"
$ws.Range("C1").Value = "MARCA"
$ws.Range("D1").Value = "DESCRIPCION"
$ws.Range("E1").Value = "CATEGORIA"
$ws.Range("F1").Value = "PRECIO"
$ws.Range("G1").Value = "STOCK"
$ws.Range("H1").Value = "FECHA_PUBLICACION"

# ---- Row 2 ----
$ws.Range("A2").Value = 1234562789
$ws.Range("C2").Value = "Fake Branch"
$ws.Range("D2").Value = "A fake description"
$ws.Range("E2").Value = "Categoria uno"
$ws.Range("F2").Value = 999999
$ws.Range("G2").Value = 1
# H2 carries a date number format (mm/dd/yy); the corrected file stores a
# literal text value there instead of a date serial, so force the cell to
# Text before writing it and then restore its original date format/style
# (copied from the still-untouched H6) without disturbing the value type.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "11/11/2020"

# ---- Row 3 ----
$ws.Range("A3").Value = 12345627289
$ws.Range("B3").Value = "Fake Name 2"
$ws.Range("D3").Value = "A fake description"
$ws.Range("E3").Value = "Categoria dos"
$ws.Range("F3").Value = 999999
$ws.Range("G3").Value = 1
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "11/11/2020"

# ---- Row 4 ----
$ws.Range("A4").Value = 1234627289
$ws.Range("B4").Value = "Fake Name 3"
$ws.Range("C4").Value = "Fake Branch 3"
$ws.Range("E4").Value = "Categoria tres"
$ws.Range("F4").Value = 999999
$ws.Range("G4").Value = 1
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "11/11/2020"

# ---- Row 5 ----
$ws.Range("A5").Value = 1246272809
$ws.Range("B5").Value = "Fake Name 4"
$ws.Range("C5").Value = "Fake Branch 4"
$ws.Range("D5").Value = "A fake description"
$ws.Range("F5").Value = 999999
$ws.Range("G5").Value = 1
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "11/11/2020"

# ---- Row 6 ----
$ws.Range("A6").Value = 12462728091
$ws.Range("B6").Value = "Fake Name 5"
$ws.Range("C6").Value = "Fake Branch 5"
$ws.Range("D6").Value = "A fake description"
$ws.Range("E6").Value = "Categoria tres"
# H6 stays empty (no value), same as before.

# Re-apply the original date cell format/style (mm/dd/yy, style index 3) to
# H2:H5 by copying formats only from the untouched, still-correctly-styled
# H6 cell, so the "@"-formatting detour above leaves no trace on these
# cells' own style while keeping their values as text.
$ws.Range("H6").Copy()
$ws.Range("H2:H5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F11").Select()
